$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.185.73"
$ws.Range("E2").Value = "  +4.76%  "
$ws.Range("D3").Value = "2.267.29"
$ws.Range("E3").Value = "  +3.57%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").Value = "254.13"
$ws.Range("E5").Value = "  -0.75%  "
$ws.Range("E6").Value = "  +2.13%  "
$ws.Range("D7").Value = "72.19"
$ws.Range("E7").Value = "  +5.33%  "
$ws.Range("D8").Value = "0.677"
$ws.Range("E8").Value = "  +19.24%  "
$ws.Range("E9").Value = "  -0.02%  "
$ws.Range("D10").Value = "40.22"
$ws.Range("E10").Value = "  +8.15%  "
$ws.Range("D11").Value = "0.0977"
$ws.Range("E11").Value = "  +4.21%  "
$ws.Range("D12").Value = "59.21"
$ws.Range("E12").Value = "  +0.60%  "
$ws.Range("D13").Value = "7.52"
$ws.Range("E13").Value = "  +6.95%  "
$ws.Range("D14").Value = "0.105"
$ws.Range("E14").Value = "  +0.70%  "
$ws.Range("D15").Value = "2.608.75"
$ws.Range("E15").Value = "  +3.95%  "
$ws.Range("D16").Value = "'14.90"
$ws.Range("E16").Value = "  +3.32%  "
$ws.Range("D17").Value = "'0.890"
$ws.Range("E17").Value = "  +2.25%  "
$ws.Range("D18").Value = "2.266.43"
$ws.Range("E18").Value = "  +4.28%  "
$ws.Range("D19").Value = "43.096.70"
$ws.Range("E19").Value = "  +4.52%  "
$ws.Range("D20").Value = "0.0₃0984"
$ws.Range("E20").Value = "  +2.76%  "
$ws.Range("D21").Value = "6.29"
$ws.Range("E21").Value = "  +1.93%  "
$ws.Range("D22").Value = "73.47"
$ws.Range("E22").Value = "  +1.95%  "
$ws.Range("D23").Value = "237.49"
$ws.Range("E23").Value = "  +1.87%  "
$ws.Range("E24").Value = "  +4.02%  "
$ws.Range("E25").Value = "  +1.95%  "
$ws.Range("D26").Value = "11.75"
$ws.Range("E26").Value = "  +0.28%  "
$ws.Range("E27").Value = "  -0.10%  "
$ws.Range("E28").Value = "  -1.67%  "
$ws.Range("E29").Value = "  +0.72%  "
$ws.Range("E30").Value = "  +2.25%  "
$ws.Range("D31").Value = "168.21"
$ws.Range("E31").Value = "  -0.59%  "
$ws.Range("D32").Value = "'21.30"
$ws.Range("E32").Value = "  +3.19%  "
$ws.Range("D33").Value = "0.128"
$ws.Range("E33").Value = "  +9.42%  "
$ws.Range("D34").Value = "6.15"
$ws.Range("E34").Value = "  +12.49%  "
$ws.Range("D35").Value = "0.0781"
$ws.Range("E35").Value = "  +4.23%  "
$ws.Range("E36").Value = "  +1.94%  "
$ws.Range("D37").Value = "28.98"
$ws.Range("E37").Value = "  +9.81%  "
$ws.Range("D38").Value = "4.74"
$ws.Range("E38").Value = "  +2.96%  "
$ws.Range("D39").Value = "4.17"
$ws.Range("E39").Value = "  -0.09%  "
$ws.Range("D40").Value = "0.0323"
$ws.Range("E40").Value = "  +8.62%  "
$ws.Range("E41").Value = "  +4.52%  "
$ws.Range("D42").Value = "5.92"
$ws.Range("E42").Value = "  +4.59%  "
$ws.Range("D43").Value = "12.58"
$ws.Range("E43").Value = "  +2.71%  "
$ws.Range("D44").Value = "'64.30"
$ws.Range("E44").Value = "  +1.22%  "
$ws.Range("D45").Value = "4.95"
$ws.Range("E45").Value = "  -0.49%  "
$ws.Range("D46").Value = "0.202"
$ws.Range("E46").Value = "  +1.91%  "
$ws.Range("D47").Value = "8.99"
$ws.Range("E47").Value = "  +4.41%  "
$ws.Range("E48").Value = "  +2.76%  "
$ws.Range("E49").Value = "  -0.17%  "
$ws.Range("D50").Value = "0.999"
$ws.Range("E50").Value = "  -0.42%  "
$ws.Range("B51").Value = "TrustWalletToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D51").Value = "'1.20"
$ws.Range("E51").Value = "  +2.54%  "
